$wb = $excel.ActiveWorkbook

# Overview sheet: mark the cbf32fc4 file (row 3) as handed back for both locales
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: mark row 3 (cbf32fc4 file) handed back, with new handback datetime
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("B3").Value = "Handed back: in sync with en-US"
$ws2.Range("G3").Value = "2016-03-08 06:24:47"

# de-de sheet: mark row 3 (cbf32fc4 file) handed back, with new handback datetime
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("B3").Value = "Handed back: in sync with en-US"
$ws3.Range("G3").Value = "2016-03-08 06:24:52"
